$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.721.11"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3
$ws.Range("D3").Value = "3.530.37"
$ws.Range("E3").Value = "  +1.12%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'606.33"
$ws.Range("E5").Value = "  -0.10%  "

# Row 6
$ws.Range("D6").Value = "'194.06"
$ws.Range("E6").Value = "  +0.63%  "

# Row 7
$ws.Range("E7").Value = "  -0.57%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.202"
$ws.Range("E9").Value = "  -4.97%  "

# Row 10
$ws.Range("E10").Value = "  -2.06%  "

# Row 11
$ws.Range("D11").Value = "'53.15"
$ws.Range("E11").Value = "  -0.30%  "

# Row 12
$ws.Range("D12").Value = "'0.0000303"
$ws.Range("E12").Value = "  -1.23%  "

# Row 13
$ws.Range("E13").Value = "  -1.41%  "

# Row 14
$ws.Range("D14").Value = "4.093.70"
$ws.Range("E14").Value = "  +1.06%  "

# Row 15
$ws.Range("D15").Value = "'591.99"
$ws.Range("E15").Value = "  -2.64%  "

# Row 16
$ws.Range("D16").Value = "'12.79"
$ws.Range("E16").Value = "  +1.24%  "

# Row 17
$ws.Range("D17").Value = "69.842.95"
$ws.Range("E17").Value = "  +0.08%  "

# Row 18
$ws.Range("D18").Value = "'18.95"
$ws.Range("E18").Value = "  +0.64%  "

# Row 19
$ws.Range("D19").Value = "3.532.73"
$ws.Range("E19").Value = "  +1.10%  "

# Row 20
$ws.Range("E20").Value = "  +1.72%  "

# Row 22
$ws.Range("D22").Value = "'17.64"
$ws.Range("E22").Value = "  -0.89%  "

# Row 23
$ws.Range("D23").Value = "'103.09"
$ws.Range("E23").Value = "  -2.22%  "

# Row 24
$ws.Range("E24").Value = "  +0.80%  "

# Row 25
$ws.Range("E25").Value = "  -0.04%  "

# Row 26
$ws.Range("D26").Value = "'3.02"
$ws.Range("E26").Value = "  -0.88%  "

# Row 27
$ws.Range("D27").Value = "'10.71"
$ws.Range("E27").Value = "  -2.16%  "

# Row 28
$ws.Range("E28").Value = "  -3.80%  "

# Row 29
$ws.Range("D29").Value = "'33.05"
$ws.Range("E29").Value = "  -3.09%  "

# Row 30
$ws.Range("B30").Value = "dogwifhat"
$ws.Range("C30").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D30").Value = "'4.26"
$ws.Range("E30").Value = "  -1.58%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.00"
$ws.Range("E31").Value = "  -1.97%  "

# Row 32
$ws.Range("D32").Value = "'12.26"
$ws.Range("E32").Value = "  -3.13%  "

# Row 33
$ws.Range("D33").Value = "'0.115"
$ws.Range("E33").Value = "  +0.06%  "

# Row 34
$ws.Range("D34").Value = "'63.29"
$ws.Range("E34").Value = "  -1.37%  "

# Row 35
$ws.Range("D35").Value = "3.818.95"
$ws.Range("E35").Value = "  +2.75%  "

# Row 36
$ws.Range("D36").Value = "'3.19"
$ws.Range("E36").Value = "  +4.80%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0818"
$ws.Range("E37").Value = "  +3.20%  "

# Row 38
$ws.Range("E38").Value = "  +0.23%  "

# Row 39
$ws.Range("D39").Value = "'512.93"
$ws.Range("E39").Value = "  -1.22%  "

# Row 40
$ws.Range("D40").Value = "'0.389"
$ws.Range("E40").Value = "  -0.18%  "

# Row 41
$ws.Range("D41").Value = "'3.55"
$ws.Range("E41").Value = "  -1.20%  "

# Row 42
$ws.Range("D42").Value = "'36.37"
$ws.Range("E42").Value = "  -0.41%  "

# Row 43
$ws.Range("E43").Value = "  -2.38%  "

# Row 44
$ws.Range("E44").Value = "  -3.23%  "

# Row 45
$ws.Range("D45").Value = "'3.34"
$ws.Range("E45").Value = "  +0.73%  "

# Row 46
$ws.Range("E46").Value = "  -1.12%  "

# Row 47
$ws.Range("E47").Value = "  -1.69%  "

# Row 48
$ws.Range("E48").Value = "  +0.16%  "

# Row 49
$ws.Range("E49").Value = "  -2.94%  "

# Row 50
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000245"
$ws.Range("E50").Value = "  +3.12%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'1.32"
$ws.Range("E51").Value = "  +1.76%  "
